$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 494, shifting the existing rows 494:534 down to 495:535.
$ws.Rows("494:494").Insert()

# Populate the newly inserted row 494 with the new record
# (the rest of the data below already shifted down automatically).
$ws.Range("A494").Value = 3
$ws.Range("B494").Value = "Femacal de La Calera"
$ws.Range("C494").Value = "Coquimbo"
$ws.Range("D494").Value = (Get-Date -Year 2022 -Month 8 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E494").Value = 5
$ws.Range("F494").Value = 100112037
$ws.Range("G494").Value = "Cebollín"
$ws.Range("H494").Value = "Sin especificar"
$ws.Range("I494").Value = "Primera"
$ws.Range("J494").Value = 110
$ws.Range("K494").Value = 7000
$ws.Range("L494").Value = 7000
$ws.Range("M494").Value = 7000
$ws.Range("N494").Value = "`$/paquete 36 unidades"
$ws.Range("O494").Value = "Provincia de Quillota"
$ws.Range("P494").Value = 194
$ws.Range("Q494").Value = 36
$ws.Range("R494").Value = "Hortaliza"
